$d = $word.ActiveDocument

# Locate the empty "List Paragraph" bullet item that immediately follows the
# "How to display image with the image path in an object" bullet. That is
# the paragraph that gets the new "Can we use interface..." text.
$count = $d.Paragraphs.Count
$anchorIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    $t2 = $t.TrimEnd([char]13, [char]7)
    if ($t2 -eq "How to display image with the image path in an object") {
        $anchorIdx = $i
        break
    }
}

$interfaceIdx = $anchorIdx + 1
$debuggingIdx = $anchorIdx + 2

# --- Empty "List Paragraph" bullet (numId=3) right after the anchor ---
# add the text "Can we use interface as a model. Do we need to create constructor in interface"
$pInterface = $d.Paragraphs.Item($interfaceIdx)
$pInterface.Range.Text = "Can we use interface as a model. Do we need to create constructor in interface"
$pInterface.Range.LanguageID = "en-US"

# --- Next empty "Normal" paragraph ---
# add the text "Debugging"
$pDebugging = $d.Paragraphs.Item($debuggingIdx)
$pDebugging.Range.Text = "Debugging"
$pDebugging.Range.LanguageID = "en-US"

# --- Insert a brand-new empty bullet-list paragraph (numId=3) right after
#     the "Debugging" paragraph, reusing the first paragraph's list
#     formatting by copy/pasting its paragraph mark (keeps the same numId
#     instead of minting a new numbering definition), then clearing its
#     text so the new paragraph stays empty. ---
$pInterface.Range.Copy()

$insertBeforeIdx = $debuggingIdx + 1
$pInsertBefore = $d.Paragraphs.Item($insertBeforeIdx)
$pasteRange = $d.Range($pInsertBefore.Range.Start, $pInsertBefore.Range.Start)
$pasteRange.Paste()

$newPara = $d.Paragraphs.Item($insertBeforeIdx)
$clearRange = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$clearRange.Delete()
